# Add speaker notes to slides 2-10 and remove the empty "Subtitle 2"
# placeholder left over on slide 7 ("Demonstration").

$p = $ppt.ActivePresentation

# ppPlaceholderBody = 2. On a notes page only the notes-body placeholder
# can be created/written through the object model, so we add it and set
# its text via TextFrame.TextRange, exactly as this host expects.
$ppPlaceholderBody = 2

$notes = @{
    2  = "Matthew"
    3  = "All of us will introduce ourselves with our Name, School, Grade, and class we were in (java class or html class)"
    4  = "Cate"
    5  = "Jared"
    6  = "Ben and Jared"
    7  = "At this point, we will give a demonstration. All of us will chirp in on this, we will discuss Saturday morning. Whatever you worked on though, you will demonstrate."
    8  = "Matthew and Cate"
    9  = "Ben"
    10 = "Questions?"
}

for ($i = 2; $i -le 10; $i++) {
    $slide = $p.Slides.Item($i)
    $notesPage = $slide.NotesPage
    $notesShape = $notesPage.Shapes.AddPlaceholder($ppPlaceholderBody)
    $notesShape.TextFrame.TextRange.Text = $notes[$i]
}

# Slide 7 ("Demonstration") still has an empty "Subtitle 2" placeholder
# shape left behind from the layout; remove it.
$slide7 = $p.Slides.Item(7)
for ($i = $slide7.Shapes.Count; $i -ge 1; $i--) {
    $shape = $slide7.Shapes.Item($i)
    if ($shape.Name -eq "Subtitle 2") {
        $shape.Delete()
    }
}
